# Auto-generated PowerShell COM-interop script
# Applies the MarkalarMallar.xlsx edit described by the commit diff:
#  1) Updates ~77 existing price cells (column C) across the Data sheet
#  2) Appends 44 new product rows (rows 2465-2508) for a new 'Bolt, Qayka, Sayba' line
#  3) Grows the AutoFilter range and the _FilterDatabase defined name to A1:H2508
#  4) Moves the saved selection to B6

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Price corrections on existing rows (column C) ---
$ws.Cells.Item(450, 3).Value = 19.899999999999999
$ws.Cells.Item(455, 3).Value = 21
$ws.Cells.Item(474, 3).Value = 16.5
$ws.Cells.Item(480, 3).Value = 11.5
$ws.Cells.Item(484, 3).Value = 7.9
$ws.Cells.Item(487, 3).Value = 0.5
$ws.Cells.Item(1242, 3).Value = 90
$ws.Cells.Item(1243, 3).Value = 105
$ws.Cells.Item(1244, 3).Value = 179
$ws.Cells.Item(1245, 3).Value = 172
$ws.Cells.Item(1306, 3).Value = 4.3
$ws.Cells.Item(1321, 3).Value = 5.6
$ws.Cells.Item(1342, 3).Value = 7
$ws.Cells.Item(1343, 3).Value = 8.6
$ws.Cells.Item(1347, 3).Value = 7.9
$ws.Cells.Item(1352, 3).Value = 5.9
$ws.Cells.Item(1482, 3).Value = 2.9
$ws.Cells.Item(1483, 3).Value = 3.2
$ws.Cells.Item(1484, 3).Value = 6.2
$ws.Cells.Item(1485, 3).Value = 11.9
$ws.Cells.Item(1626, 3).Value = 0.7
$ws.Cells.Item(1627, 3).Value = 0.7
$ws.Cells.Item(1628, 3).Value = 0.7
$ws.Cells.Item(1629, 3).Value = 0.7
$ws.Cells.Item(1630, 3).Value = 0.7
$ws.Cells.Item(1631, 3).Value = 0.7
$ws.Cells.Item(1632, 3).Value = 0.7
$ws.Cells.Item(1633, 3).Value = 0.7
$ws.Cells.Item(1634, 3).Value = 0.7
$ws.Cells.Item(1635, 3).Value = 0.7
$ws.Cells.Item(1636, 3).Value = 0.7
$ws.Cells.Item(1637, 3).Value = 0.7
$ws.Cells.Item(1638, 3).Value = 0.7
$ws.Cells.Item(1639, 3).Value = 0.7
$ws.Cells.Item(1641, 3).Value = 0.7
$ws.Cells.Item(1648, 3).Value = 15
$ws.Cells.Item(1727, 3).Value = 1.7
$ws.Cells.Item(1728, 3).Value = 1.7
$ws.Cells.Item(1729, 3).Value = 1.7
$ws.Cells.Item(1730, 3).Value = 1.5
$ws.Cells.Item(1731, 3).Value = 1.5
$ws.Cells.Item(1732, 3).Value = 5.8
$ws.Cells.Item(1744, 3).Value = 97
$ws.Cells.Item(1793, 3).Value = 38
$ws.Cells.Item(1848, 3).Value = 43
$ws.Cells.Item(1851, 3).Value = 75
$ws.Cells.Item(1852, 3).Value = 33.200000000000003
$ws.Cells.Item(1853, 3).Value = 51
$ws.Cells.Item(1873, 3).Value = 26.2
$ws.Cells.Item(1878, 3).Value = 629
$ws.Cells.Item(1881, 3).Value = 219
$ws.Cells.Item(1884, 3).Value = 43
$ws.Cells.Item(1888, 3).Value = 121
$ws.Cells.Item(1893, 3).Value = 95
$ws.Cells.Item(1901, 3).Value = 46.9
$ws.Cells.Item(1907, 3).Value = 73
$ws.Cells.Item(1909, 3).Value = 95
$ws.Cells.Item(1910, 3).Value = 106
$ws.Cells.Item(1917, 3).Value = 129
$ws.Cells.Item(2059, 3).Value = 0.72
$ws.Cells.Item(2066, 3).Value = 0.72
$ws.Cells.Item(2072, 3).Value = 1.07
$ws.Cells.Item(2155, 3).Value = 1.6
$ws.Cells.Item(2170, 3).Value = 4.8
$ws.Cells.Item(2284, 3).Value = 2.5499999999999998
$ws.Cells.Item(2293, 3).Value = 3.5
$ws.Cells.Item(2320, 3).Value = 0.37
$ws.Cells.Item(2327, 3).Value = 0.17
$ws.Cells.Item(2330, 3).Value = 0.85
$ws.Cells.Item(2335, 3).Value = 2.1
$ws.Cells.Item(2341, 3).Value = 5.9
$ws.Cells.Item(2353, 3).Value = 0.23
$ws.Cells.Item(2442, 3).Value = 7.9
$ws.Cells.Item(2443, 3).Value = 8.3000000000000007
$ws.Cells.Item(2448, 3).Value = 7.9
$ws.Cells.Item(2449, 3).Value = 7.9
$ws.Cells.Item(2450, 3).Value = 8.3000000000000007

# --- 2) Append new rows 2465-2508 ---
$ws.Cells.Item(2465, 1).Value = "TM.241105101474"
$ws.Cells.Item(2465, 2).Value = "QAYKA M6 BQM6 Q/25"
$ws.Cells.Item(2465, 3).Value = 4.3
$ws.Cells.Item(2465, 4).Value = "Xırdavat"
$ws.Cells.Item(2465, 5).Value = "Xırdavat və əl alətləri"
$ws.Cells.Item(2465, 6).Value = "Xırdavat məhsulları"
$ws.Cells.Item(2465, 7).Value = "Bolt, Qayka, Şayba"
$ws.Cells.Item(2465, 8).Value = "YIWU HAOXING"
$ws.Cells.Item(2466, 1).Value = "TM.241105101475"
$ws.Cells.Item(2466, 2).Value = "QAYKA M8 BQM8 Q/25"
$ws.Cells.Item(2466, 3).Value = 4
$ws.Cells.Item(2466, 4).Value = "Xırdavat"
$ws.Cells.Item(2466, 5).Value = "Xırdavat və əl alətləri"
$ws.Cells.Item(2466, 6).Value = "Xırdavat məhsulları"
$ws.Cells.Item(2466, 7).Value = "Bolt, Qayka, Şayba"
$ws.Cells.Item(2466, 8).Value = "YIWU HAOXING"
$ws.Cells.Item(2467, 1).Value = "TM.241105101476"
$ws.Cells.Item(2467, 2).Value = "QAYKA M10 BQM10 Q/25"
$ws.Cells.Item(2467, 3).Value = 3.9
$ws.Cells.Item(2467, 4).Value = "Xırdavat"
$ws.Cells.Item(2467, 5).Value = "Xırdavat və əl alətləri"
$ws.Cells.Item(2467, 6).Value = "Xırdavat məhsulları"
$ws.Cells.Item(2467, 7).Value = "Bolt, Qayka, Şayba"
$ws.Cells.Item(2467, 8).Value = "YIWU HAOXING"
$ws.Cells.Item(2468, 1).Value = "TM.241105101477"
$ws.Cells.Item(2468, 2).Value = "BOLT M6 15MM BM6-15MM Q/25"
$ws.Cells.Item(2468, 3).Value = 4.3
$ws.Cells.Item(2468, 4).Value = "Xırdavat"
$ws.Cells.Item(2468, 5).Value = "Xırdavat və əl alətləri"
$ws.Cells.Item(2468, 6).Value = "Xırdavat məhsulları"
$ws.Cells.Item(2468, 7).Value = "Bolt, Qayka, Şayba"
$ws.Cells.Item(2468, 8).Value = "YIWU HAOXING"
$ws.Cells.Item(2469, 1).Value = "TM.241105101478"
$ws.Cells.Item(2469, 2).Value = "BOLT M6 20MM BM6-20MM Q/25"
$ws.Cells.Item(2469, 3).Value = 4.2
$ws.Cells.Item(2469, 4).Value = "Xırdavat"
$ws.Cells.Item(2469, 5).Value = "Xırdavat və əl alətləri"
$ws.Cells.Item(2469, 6).Value = "Xırdavat məhsulları"
$ws.Cells.Item(2469, 7).Value = "Bolt, Qayka, Şayba"
$ws.Cells.Item(2469, 8).Value = "YIWU HAOXING"
$ws.Cells.Item(2470, 1).Value = "TM.241105101479"
$ws.Cells.Item(2470, 2).Value = "BOLT M6 25MM BM6-25MM Q/25"
$ws.Cells.Item(2470, 3).Value = 4.2
$ws.Cells.Item(2470, 4).Value = "Xırdavat"
$ws.Cells.Item(2470, 5).Value = "Xırdavat və əl alətləri"
$ws.Cells.Item(2470, 6).Value = "Xırdavat məhsulları"
$ws.Cells.Item(2470, 7).Value = "Bolt, Qayka, Şayba"
$ws.Cells.Item(2470, 8).Value = "YIWU HAOXING"
$ws.Cells.Item(2471, 1).Value = "TM.241105101480"
$ws.Cells.Item(2471, 2).Value = "BOLT M6 30MM BM6-30MM Q/25"
$ws.Cells.Item(2471, 3).Value = 4.0999999999999996
$ws.Cells.Item(2471, 4).Value = "Xırdavat"
$ws.Cells.Item(2471, 5).Value = "Xırdavat və əl alətləri"
$ws.Cells.Item(2471, 6).Value = "Xırdavat məhsulları"
$ws.Cells.Item(2471, 7).Value = "Bolt, Qayka, Şayba"
$ws.Cells.Item(2471, 8).Value = "YIWU HAOXING"
$ws.Cells.Item(2472, 1).Value = "TM.241105101481"
$ws.Cells.Item(2472, 2).Value = "BOLT M6 40MM BM6-40MM Q/25"
$ws.Cells.Item(2472, 3).Value = 4.0999999999999996
$ws.Cells.Item(2472, 4).Value = "Xırdavat"
$ws.Cells.Item(2472, 5).Value = "Xırdavat və əl alətləri"
$ws.Cells.Item(2472, 6).Value = "Xırdavat məhsulları"
$ws.Cells.Item(2472, 7).Value = "Bolt, Qayka, Şayba"
$ws.Cells.Item(2472, 8).Value = "YIWU HAOXING"
$ws.Cells.Item(2473, 1).Value = "TM.241105101482"
$ws.Cells.Item(2473, 2).Value = "BOLT M6 50MM BM6-50MM Q/25"
$ws.Cells.Item(2473, 3).Value = 4.0999999999999996
$ws.Cells.Item(2473, 4).Value = "Xırdavat"
$ws.Cells.Item(2473, 5).Value = "Xırdavat və əl alətləri"
$ws.Cells.Item(2473, 6).Value = "Xırdavat məhsulları"
$ws.Cells.Item(2473, 7).Value = "Bolt, Qayka, Şayba"
$ws.Cells.Item(2473, 8).Value = "YIWU HAOXING"
$ws.Cells.Item(2474, 1).Value = "TM.241105101483"
$ws.Cells.Item(2474, 2).Value = "BOLT M8 15MM BM8-15MM Q/25"
$ws.Cells.Item(2474, 3).Value = 4
$ws.Cells.Item(2474, 4).Value = "Xırdavat"
$ws.Cells.Item(2474, 5).Value = "Xırdavat və əl alətləri"
$ws.Cells.Item(2474, 6).Value = "Xırdavat məhsulları"
$ws.Cells.Item(2474, 7).Value = "Bolt, Qayka, Şayba"
$ws.Cells.Item(2474, 8).Value = "YIWU HAOXING"
$ws.Cells.Item(2475, 1).Value = "TM.241105101484"
$ws.Cells.Item(2475, 2).Value = "BOLT M8 20MM BM8-20MM Q/25"
$ws.Cells.Item(2475, 3).Value = 4
$ws.Cells.Item(2475, 4).Value = "Xırdavat"
$ws.Cells.Item(2475, 5).Value = "Xırdavat və əl alətləri"
$ws.Cells.Item(2475, 6).Value = "Xırdavat məhsulları"
$ws.Cells.Item(2475, 7).Value = "Bolt, Qayka, Şayba"
$ws.Cells.Item(2475, 8).Value = "YIWU HAOXING"
$ws.Cells.Item(2476, 1).Value = "TM.241105101485"
$ws.Cells.Item(2476, 2).Value = "BOLT M8 25MM BM8-25MM Q/25"
$ws.Cells.Item(2476, 3).Value = 3.9
$ws.Cells.Item(2476, 4).Value = "Xırdavat"
$ws.Cells.Item(2476, 5).Value = "Xırdavat və əl alətləri"
$ws.Cells.Item(2476, 6).Value = "Xırdavat məhsulları"
$ws.Cells.Item(2476, 7).Value = "Bolt, Qayka, Şayba"
$ws.Cells.Item(2476, 8).Value = "YIWU HAOXING"
$ws.Cells.Item(2477, 1).Value = "TM.241105101486"
$ws.Cells.Item(2477, 2).Value = "BOLT M8 30MM BM8-30MM Q/25"
$ws.Cells.Item(2477, 3).Value = 3.9
$ws.Cells.Item(2477, 4).Value = "Xırdavat"
$ws.Cells.Item(2477, 5).Value = "Xırdavat və əl alətləri"
$ws.Cells.Item(2477, 6).Value = "Xırdavat məhsulları"
$ws.Cells.Item(2477, 7).Value = "Bolt, Qayka, Şayba"
$ws.Cells.Item(2477, 8).Value = "YIWU HAOXING"
$ws.Cells.Item(2478, 1).Value = "TM.241105101487"
$ws.Cells.Item(2478, 2).Value = "BOLT M8 40MM BM8-40MM Q/25"
$ws.Cells.Item(2478, 3).Value = 3.9
$ws.Cells.Item(2478, 4).Value = "Xırdavat"
$ws.Cells.Item(2478, 5).Value = "Xırdavat və əl alətləri"
$ws.Cells.Item(2478, 6).Value = "Xırdavat məhsulları"
$ws.Cells.Item(2478, 7).Value = "Bolt, Qayka, Şayba"
$ws.Cells.Item(2478, 8).Value = "YIWU HAOXING"
$ws.Cells.Item(2479, 1).Value = "TM.241105101488"
$ws.Cells.Item(2479, 2).Value = "BOLT M8 50MM BM8-50MM Q/25"
$ws.Cells.Item(2479, 3).Value = 3.9
$ws.Cells.Item(2479, 4).Value = "Xırdavat"
$ws.Cells.Item(2479, 5).Value = "Xırdavat və əl alətləri"
$ws.Cells.Item(2479, 6).Value = "Xırdavat məhsulları"
$ws.Cells.Item(2479, 7).Value = "Bolt, Qayka, Şayba"
$ws.Cells.Item(2479, 8).Value = "YIWU HAOXING"
$ws.Cells.Item(2480, 1).Value = "TM.241105101489"
$ws.Cells.Item(2480, 2).Value = "BOLT M8 60MM BM8-60MM Q/25"
$ws.Cells.Item(2480, 3).Value = 3.9
$ws.Cells.Item(2480, 4).Value = "Xırdavat"
$ws.Cells.Item(2480, 5).Value = "Xırdavat və əl alətləri"
$ws.Cells.Item(2480, 6).Value = "Xırdavat məhsulları"
$ws.Cells.Item(2480, 7).Value = "Bolt, Qayka, Şayba"
$ws.Cells.Item(2480, 8).Value = "YIWU HAOXING"
$ws.Cells.Item(2481, 1).Value = "TM.241105101490"
$ws.Cells.Item(2481, 2).Value = "BOLT M8 70MM BM8-70MM Q/25"
$ws.Cells.Item(2481, 3).Value = 3.9
$ws.Cells.Item(2481, 4).Value = "Xırdavat"
$ws.Cells.Item(2481, 5).Value = "Xırdavat və əl alətləri"
$ws.Cells.Item(2481, 6).Value = "Xırdavat məhsulları"
$ws.Cells.Item(2481, 7).Value = "Bolt, Qayka, Şayba"
$ws.Cells.Item(2481, 8).Value = "YIWU HAOXING"
$ws.Cells.Item(2482, 1).Value = "TM.241105101491"
$ws.Cells.Item(2482, 2).Value = "BOLT M8 80MM BM8-80MM Q/25"
$ws.Cells.Item(2482, 3).Value = 3.9
$ws.Cells.Item(2482, 4).Value = "Xırdavat"
$ws.Cells.Item(2482, 5).Value = "Xırdavat və əl alətləri"
$ws.Cells.Item(2482, 6).Value = "Xırdavat məhsulları"
$ws.Cells.Item(2482, 7).Value = "Bolt, Qayka, Şayba"
$ws.Cells.Item(2482, 8).Value = "YIWU HAOXING"
$ws.Cells.Item(2483, 1).Value = "TM.241105101492"
$ws.Cells.Item(2483, 2).Value = "BOLT M8 100MM BM8-100MM Q/25"
$ws.Cells.Item(2483, 3).Value = 3.9
$ws.Cells.Item(2483, 4).Value = "Xırdavat"
$ws.Cells.Item(2483, 5).Value = "Xırdavat və əl alətləri"
$ws.Cells.Item(2483, 6).Value = "Xırdavat məhsulları"
$ws.Cells.Item(2483, 7).Value = "Bolt, Qayka, Şayba"
$ws.Cells.Item(2483, 8).Value = "YIWU HAOXING"
$ws.Cells.Item(2484, 1).Value = "TM.241105101493"
$ws.Cells.Item(2484, 2).Value = "BOLT M8 120MM BM8-120MM Q/25"
$ws.Cells.Item(2484, 3).Value = 4
$ws.Cells.Item(2484, 4).Value = "Xırdavat"
$ws.Cells.Item(2484, 5).Value = "Xırdavat və əl alətləri"
$ws.Cells.Item(2484, 6).Value = "Xırdavat məhsulları"
$ws.Cells.Item(2484, 7).Value = "Bolt, Qayka, Şayba"
$ws.Cells.Item(2484, 8).Value = "YIWU HAOXING"
$ws.Cells.Item(2485, 1).Value = "TM.241105101494"
$ws.Cells.Item(2485, 2).Value = "BOLT M8 150MM BM8-150MM Q/25"
$ws.Cells.Item(2485, 3).Value = 4
$ws.Cells.Item(2485, 4).Value = "Xırdavat"
$ws.Cells.Item(2485, 5).Value = "Xırdavat və əl alətləri"
$ws.Cells.Item(2485, 6).Value = "Xırdavat məhsulları"
$ws.Cells.Item(2485, 7).Value = "Bolt, Qayka, Şayba"
$ws.Cells.Item(2485, 8).Value = "YIWU HAOXING"
$ws.Cells.Item(2486, 1).Value = "TM.241105101495"
$ws.Cells.Item(2486, 2).Value = "BOLT M10 20MM BM10-20MM Q/25"
$ws.Cells.Item(2486, 3).Value = 3.8
$ws.Cells.Item(2486, 4).Value = "Xırdavat"
$ws.Cells.Item(2486, 5).Value = "Xırdavat və əl alətləri"
$ws.Cells.Item(2486, 6).Value = "Xırdavat məhsulları"
$ws.Cells.Item(2486, 7).Value = "Bolt, Qayka, Şayba"
$ws.Cells.Item(2486, 8).Value = "YIWU HAOXING"
$ws.Cells.Item(2487, 1).Value = "TM.241105101496"
$ws.Cells.Item(2487, 2).Value = "BOLT M10 25MM BM10-25MM Q/25"
$ws.Cells.Item(2487, 3).Value = 3.8
$ws.Cells.Item(2487, 4).Value = "Xırdavat"
$ws.Cells.Item(2487, 5).Value = "Xırdavat və əl alətləri"
$ws.Cells.Item(2487, 6).Value = "Xırdavat məhsulları"
$ws.Cells.Item(2487, 7).Value = "Bolt, Qayka, Şayba"
$ws.Cells.Item(2487, 8).Value = "YIWU HAOXING"
$ws.Cells.Item(2488, 1).Value = "TM.241105101497"
$ws.Cells.Item(2488, 2).Value = "BOLT M10 30MM BM10-30MM Q/25"
$ws.Cells.Item(2488, 3).Value = 3.8
$ws.Cells.Item(2488, 4).Value = "Xırdavat"
$ws.Cells.Item(2488, 5).Value = "Xırdavat və əl alətləri"
$ws.Cells.Item(2488, 6).Value = "Xırdavat məhsulları"
$ws.Cells.Item(2488, 7).Value = "Bolt, Qayka, Şayba"
$ws.Cells.Item(2488, 8).Value = "YIWU HAOXING"
$ws.Cells.Item(2489, 1).Value = "TM.241105101498"
$ws.Cells.Item(2489, 2).Value = "BOLT M10 40MM BM10-40MM Q/25"
$ws.Cells.Item(2489, 3).Value = 3.8
$ws.Cells.Item(2489, 4).Value = "Xırdavat"
$ws.Cells.Item(2489, 5).Value = "Xırdavat və əl alətləri"
$ws.Cells.Item(2489, 6).Value = "Xırdavat məhsulları"
$ws.Cells.Item(2489, 7).Value = "Bolt, Qayka, Şayba"
$ws.Cells.Item(2489, 8).Value = "YIWU HAOXING"
$ws.Cells.Item(2490, 1).Value = "TM.241105101499"
$ws.Cells.Item(2490, 2).Value = "BOLT M10 50MM BM10-50MM Q/25"
$ws.Cells.Item(2490, 3).Value = 3.8
$ws.Cells.Item(2490, 4).Value = "Xırdavat"
$ws.Cells.Item(2490, 5).Value = "Xırdavat və əl alətləri"
$ws.Cells.Item(2490, 6).Value = "Xırdavat məhsulları"
$ws.Cells.Item(2490, 7).Value = "Bolt, Qayka, Şayba"
$ws.Cells.Item(2490, 8).Value = "YIWU HAOXING"
$ws.Cells.Item(2491, 1).Value = "TM.241105101500"
$ws.Cells.Item(2491, 2).Value = "BOLT M10 60MM BM10-60MM Q/25"
$ws.Cells.Item(2491, 3).Value = 3.8
$ws.Cells.Item(2491, 4).Value = "Xırdavat"
$ws.Cells.Item(2491, 5).Value = "Xırdavat və əl alətləri"
$ws.Cells.Item(2491, 6).Value = "Xırdavat məhsulları"
$ws.Cells.Item(2491, 7).Value = "Bolt, Qayka, Şayba"
$ws.Cells.Item(2491, 8).Value = "YIWU HAOXING"
$ws.Cells.Item(2492, 1).Value = "TM.241105101501"
$ws.Cells.Item(2492, 2).Value = "BOLT M10 70MM BM10-70MM Q/25"
$ws.Cells.Item(2492, 3).Value = 3.8
$ws.Cells.Item(2492, 4).Value = "Xırdavat"
$ws.Cells.Item(2492, 5).Value = "Xırdavat və əl alətləri"
$ws.Cells.Item(2492, 6).Value = "Xırdavat məhsulları"
$ws.Cells.Item(2492, 7).Value = "Bolt, Qayka, Şayba"
$ws.Cells.Item(2492, 8).Value = "YIWU HAOXING"
$ws.Cells.Item(2493, 1).Value = "TM.241105101502"
$ws.Cells.Item(2493, 2).Value = "BOLT M10 100MM BM10-100MM Q/25"
$ws.Cells.Item(2493, 3).Value = 3.8
$ws.Cells.Item(2493, 4).Value = "Xırdavat"
$ws.Cells.Item(2493, 5).Value = "Xırdavat və əl alətləri"
$ws.Cells.Item(2493, 6).Value = "Xırdavat məhsulları"
$ws.Cells.Item(2493, 7).Value = "Bolt, Qayka, Şayba"
$ws.Cells.Item(2493, 8).Value = "YIWU HAOXING"
$ws.Cells.Item(2494, 1).Value = "TM.241105101503"
$ws.Cells.Item(2494, 2).Value = "BOLT M10 120MM BM10-120MM Q/25"
$ws.Cells.Item(2494, 3).Value = 3.8
$ws.Cells.Item(2494, 4).Value = "Xırdavat"
$ws.Cells.Item(2494, 5).Value = "Xırdavat və əl alətləri"
$ws.Cells.Item(2494, 6).Value = "Xırdavat məhsulları"
$ws.Cells.Item(2494, 7).Value = "Bolt, Qayka, Şayba"
$ws.Cells.Item(2494, 8).Value = "YIWU HAOXING"
$ws.Cells.Item(2495, 1).Value = "TM.241105101504"
$ws.Cells.Item(2495, 2).Value = "ŞAYBA 6X16 BS6*16 Q/25"
$ws.Cells.Item(2495, 3).Value = 3.9
$ws.Cells.Item(2495, 4).Value = "Xırdavat"
$ws.Cells.Item(2495, 5).Value = "Xırdavat və əl alətləri"
$ws.Cells.Item(2495, 6).Value = "Xırdavat məhsulları"
$ws.Cells.Item(2495, 7).Value = "Bolt, Qayka, Şayba"
$ws.Cells.Item(2495, 8).Value = "YIWU HAOXING"
$ws.Cells.Item(2496, 1).Value = "TM.241105101505"
$ws.Cells.Item(2496, 2).Value = "ŞAYBA 6X18 BS6*18 Q/25"
$ws.Cells.Item(2496, 3).Value = 3.9
$ws.Cells.Item(2496, 4).Value = "Xırdavat"
$ws.Cells.Item(2496, 5).Value = "Xırdavat və əl alətləri"
$ws.Cells.Item(2496, 6).Value = "Xırdavat məhsulları"
$ws.Cells.Item(2496, 7).Value = "Bolt, Qayka, Şayba"
$ws.Cells.Item(2496, 8).Value = "YIWU HAOXING"
$ws.Cells.Item(2497, 1).Value = "TM.241105101506"
$ws.Cells.Item(2497, 2).Value = "ŞAYBA 6X22 BS6*22 Q/25"
$ws.Cells.Item(2497, 3).Value = 3.9
$ws.Cells.Item(2497, 4).Value = "Xırdavat"
$ws.Cells.Item(2497, 5).Value = "Xırdavat və əl alətləri"
$ws.Cells.Item(2497, 6).Value = "Xırdavat məhsulları"
$ws.Cells.Item(2497, 7).Value = "Bolt, Qayka, Şayba"
$ws.Cells.Item(2497, 8).Value = "YIWU HAOXING"
$ws.Cells.Item(2498, 1).Value = "TM.241105101507"
$ws.Cells.Item(2498, 2).Value = "ŞAYBA 6X30 BS6*30 Q/25"
$ws.Cells.Item(2498, 3).Value = 3.9
$ws.Cells.Item(2498, 4).Value = "Xırdavat"
$ws.Cells.Item(2498, 5).Value = "Xırdavat və əl alətləri"
$ws.Cells.Item(2498, 6).Value = "Xırdavat məhsulları"
$ws.Cells.Item(2498, 7).Value = "Bolt, Qayka, Şayba"
$ws.Cells.Item(2498, 8).Value = "YIWU HAOXING"
$ws.Cells.Item(2499, 1).Value = "TM.241105101508"
$ws.Cells.Item(2499, 2).Value = "ŞAYBA 8X18 BS8*18 Q/25"
$ws.Cells.Item(2499, 3).Value = 3.9
$ws.Cells.Item(2499, 4).Value = "Xırdavat"
$ws.Cells.Item(2499, 5).Value = "Xırdavat və əl alətləri"
$ws.Cells.Item(2499, 6).Value = "Xırdavat məhsulları"
$ws.Cells.Item(2499, 7).Value = "Bolt, Qayka, Şayba"
$ws.Cells.Item(2499, 8).Value = "YIWU HAOXING"
$ws.Cells.Item(2500, 1).Value = "TM.241105101509"
$ws.Cells.Item(2500, 2).Value = "ŞAYBA 8X22 BS8*22 Q/25"
$ws.Cells.Item(2500, 3).Value = 3.9
$ws.Cells.Item(2500, 4).Value = "Xırdavat"
$ws.Cells.Item(2500, 5).Value = "Xırdavat və əl alətləri"
$ws.Cells.Item(2500, 6).Value = "Xırdavat məhsulları"
$ws.Cells.Item(2500, 7).Value = "Bolt, Qayka, Şayba"
$ws.Cells.Item(2500, 8).Value = "YIWU HAOXING"
$ws.Cells.Item(2501, 1).Value = "TM.241105101510"
$ws.Cells.Item(2501, 2).Value = "ŞAYBA 8X30 BS8*30 Q/25"
$ws.Cells.Item(2501, 3).Value = 3.9
$ws.Cells.Item(2501, 4).Value = "Xırdavat"
$ws.Cells.Item(2501, 5).Value = "Xırdavat və əl alətləri"
$ws.Cells.Item(2501, 6).Value = "Xırdavat məhsulları"
$ws.Cells.Item(2501, 7).Value = "Bolt, Qayka, Şayba"
$ws.Cells.Item(2501, 8).Value = "YIWU HAOXING"
$ws.Cells.Item(2502, 1).Value = "TM.241105101511"
$ws.Cells.Item(2502, 2).Value = "ŞAYBA 8X40 BS8*40 Q/25"
$ws.Cells.Item(2502, 3).Value = 3.9
$ws.Cells.Item(2502, 4).Value = "Xırdavat"
$ws.Cells.Item(2502, 5).Value = "Xırdavat və əl alətləri"
$ws.Cells.Item(2502, 6).Value = "Xırdavat məhsulları"
$ws.Cells.Item(2502, 7).Value = "Bolt, Qayka, Şayba"
$ws.Cells.Item(2502, 8).Value = "YIWU HAOXING"
$ws.Cells.Item(2503, 1).Value = "TM.241105101512"
$ws.Cells.Item(2503, 2).Value = "ŞAYBA 10X25 BS10*25 Q/25"
$ws.Cells.Item(2503, 3).Value = 3.9
$ws.Cells.Item(2503, 4).Value = "Xırdavat"
$ws.Cells.Item(2503, 5).Value = "Xırdavat və əl alətləri"
$ws.Cells.Item(2503, 6).Value = "Xırdavat məhsulları"
$ws.Cells.Item(2503, 7).Value = "Bolt, Qayka, Şayba"
$ws.Cells.Item(2503, 8).Value = "YIWU HAOXING"
$ws.Cells.Item(2504, 1).Value = "TM.241105101513"
$ws.Cells.Item(2504, 2).Value = "ŞAYBA 10X30 BS10*30 Q/25"
$ws.Cells.Item(2504, 3).Value = 3.9
$ws.Cells.Item(2504, 4).Value = "Xırdavat"
$ws.Cells.Item(2504, 5).Value = "Xırdavat və əl alətləri"
$ws.Cells.Item(2504, 6).Value = "Xırdavat məhsulları"
$ws.Cells.Item(2504, 7).Value = "Bolt, Qayka, Şayba"
$ws.Cells.Item(2504, 8).Value = "YIWU HAOXING"
$ws.Cells.Item(2505, 1).Value = "TM.241105091514"
$ws.Cells.Item(2505, 2).Value = "BOLT SAMAREZ 8X75MM HSW875 Q/25"
$ws.Cells.Item(2505, 3).Value = 4.2
$ws.Cells.Item(2505, 4).Value = "Xırdavat"
$ws.Cells.Item(2505, 5).Value = "Xırdavat və əl alətləri"
$ws.Cells.Item(2505, 6).Value = "Xırdavat məhsulları"
$ws.Cells.Item(2505, 7).Value = "Bolt, Qayka, Şayba"
$ws.Cells.Item(2505, 8).Value = "YIWU HAOXING"
$ws.Cells.Item(2506, 1).Value = "TM.241105091515"
$ws.Cells.Item(2506, 2).Value = "BOLT SAMAREZ 10X75MM HSW1075 Q/25"
$ws.Cells.Item(2506, 3).Value = 4.2
$ws.Cells.Item(2506, 4).Value = "Xırdavat"
$ws.Cells.Item(2506, 5).Value = "Xırdavat və əl alətləri"
$ws.Cells.Item(2506, 6).Value = "Xırdavat məhsulları"
$ws.Cells.Item(2506, 7).Value = "Bolt, Qayka, Şayba"
$ws.Cells.Item(2506, 8).Value = "YIWU HAOXING"
$ws.Cells.Item(2507, 1).Value = "TM.241105091516"
$ws.Cells.Item(2507, 2).Value = "BOLT SAMAREZ 8X100MM HSW810 Q/25"
$ws.Cells.Item(2507, 3).Value = 4.2
$ws.Cells.Item(2507, 4).Value = "Xırdavat"
$ws.Cells.Item(2507, 5).Value = "Xırdavat və əl alətləri"
$ws.Cells.Item(2507, 6).Value = "Xırdavat məhsulları"
$ws.Cells.Item(2507, 7).Value = "Bolt, Qayka, Şayba"
$ws.Cells.Item(2507, 8).Value = "YIWU HAOXING"
$ws.Cells.Item(2508, 1).Value = "TM.241105091517"
$ws.Cells.Item(2508, 2).Value = "BOLT SAMAREZ 10X100MM HSW1010 Q/25"
$ws.Cells.Item(2508, 3).Value = 4.2
$ws.Cells.Item(2508, 4).Value = "Xırdavat"
$ws.Cells.Item(2508, 5).Value = "Xırdavat və əl alətləri"
$ws.Cells.Item(2508, 6).Value = "Xırdavat məhsulları"
$ws.Cells.Item(2508, 7).Value = "Bolt, Qayka, Şayba"
$ws.Cells.Item(2508, 8).Value = "YIWU HAOXING"

# --- 3) Grow the AutoFilter range to cover the new rows ---
if ($ws.AutoFilter) {
    $ws.AutoFilterMode = $false
}
$ws.Range("A1:H2508").AutoFilter()

# --- Keep the _xlnm._FilterDatabase defined name range in sync ---
$filterName = $wb.Names.Item(1)
$filterName.RefersTo = "=Data!`$A`$1:`$H`$2508"

# --- 4) Restore the saved cursor / selection position ---
$ws.Range("B6").Select()

